$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(4, 9).Value = 'ba'
$ws.Cells.Item(4, 10).Value = 'Appreciation'
$ws.Cells.Item(16, 9).Value = 'ba'
$ws.Cells.Item(16, 10).Value = 'Appreciation'
$ws.Cells.Item(24, 9).Value = 'ba'
$ws.Cells.Item(24, 10).Value = 'Appreciation'
$ws.Cells.Item(35, 9).Value = 'sd'
$ws.Cells.Item(35, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(36, 9).Value = 'ba'
$ws.Cells.Item(36, 10).Value = 'Appreciation'
$ws.Cells.Item(47, 9).Value = 'sv'
$ws.Cells.Item(47, 10).Value = 'Statement-opinion'
$ws.Cells.Item(56, 9).Value = 'b'
$ws.Cells.Item(56, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(78, 9).Value = 'b'
$ws.Cells.Item(78, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(91, 9).Value = '%'
$ws.Cells.Item(91, 10).Value = 'Uninterpretable'
$ws.Cells.Item(99, 9).Value = 'sv'
$ws.Cells.Item(99, 10).Value = 'Statement-opinion'
$ws.Cells.Item(122, 9).Value = 'sv'
$ws.Cells.Item(122, 10).Value = 'Statement-opinion'
$ws.Cells.Item(126, 9).Value = 'b'
$ws.Cells.Item(126, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(131, 9).Value = 'ba'
$ws.Cells.Item(131, 10).Value = 'Appreciation'
$ws.Cells.Item(134, 9).Value = 'sv'
$ws.Cells.Item(134, 10).Value = 'Statement-opinion'
$ws.Cells.Item(144, 9).Value = 'ba'
$ws.Cells.Item(144, 10).Value = 'Appreciation'
$ws.Cells.Item(146, 9).Value = 'ba'
$ws.Cells.Item(146, 10).Value = 'Appreciation'
$ws.Cells.Item(149, 9).Value = 'ba'
$ws.Cells.Item(149, 10).Value = 'Appreciation'
$ws.Cells.Item(151, 9).Value = 'sv'
$ws.Cells.Item(151, 10).Value = 'Statement-opinion'
$ws.Cells.Item(152, 9).Value = 'sv'
$ws.Cells.Item(152, 10).Value = 'Statement-opinion'
$ws.Cells.Item(153, 9).Value = 'b'
$ws.Cells.Item(153, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(160, 9).Value = 'sv'
$ws.Cells.Item(160, 10).Value = 'Statement-opinion'
$ws.Cells.Item(164, 9).Value = 'sv'
$ws.Cells.Item(164, 10).Value = 'Statement-opinion'
$ws.Cells.Item(167, 9).Value = 'qy'
$ws.Cells.Item(167, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(173, 9).Value = 'aa'
$ws.Cells.Item(173, 10).Value = 'Agree/Accept'
$ws.Cells.Item(175, 9).Value = 'sv'
$ws.Cells.Item(175, 10).Value = 'Statement-opinion'
$ws.Cells.Item(178, 9).Value = 'sv'
$ws.Cells.Item(178, 10).Value = 'Statement-opinion'
$ws.Cells.Item(183, 9).Value = 'sd'
$ws.Cells.Item(183, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(187, 9).Value = 'b'
$ws.Cells.Item(187, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(189, 9).Value = 'aa'
$ws.Cells.Item(189, 10).Value = 'Agree/Accept'
$ws.Cells.Item(190, 9).Value = 'sd'
$ws.Cells.Item(190, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(193, 9).Value = 'aa'
$ws.Cells.Item(193, 10).Value = 'Agree/Accept'
$ws.Cells.Item(195, 9).Value = 'sd'
$ws.Cells.Item(195, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(204, 9).Value = 'sd'
$ws.Cells.Item(204, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(209, 9).Value = 'sd'
$ws.Cells.Item(209, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(215, 9).Value = 'sd'
$ws.Cells.Item(215, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(219, 9).Value = 'ba'
$ws.Cells.Item(219, 10).Value = 'Appreciation'
$ws.Cells.Item(221, 9).Value = 'ba'
$ws.Cells.Item(221, 10).Value = 'Appreciation'
$ws.Cells.Item(234, 9).Value = 'ba'
$ws.Cells.Item(234, 10).Value = 'Appreciation'
$ws.Cells.Item(248, 9).Value = 'ba'
$ws.Cells.Item(248, 10).Value = 'Appreciation'
$ws.Cells.Item(251, 9).Value = 'ba'
$ws.Cells.Item(251, 10).Value = 'Appreciation'
$ws.Cells.Item(252, 9).Value = 'sd'
$ws.Cells.Item(252, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(253, 9).Value = 'sv'
$ws.Cells.Item(253, 10).Value = 'Statement-opinion'
$ws.Cells.Item(268, 9).Value = 'ba'
$ws.Cells.Item(268, 10).Value = 'Appreciation'
$ws.Cells.Item(272, 9).Value = 'aa'
$ws.Cells.Item(272, 10).Value = 'Agree/Accept'
$ws.Cells.Item(275, 9).Value = 'sd'
$ws.Cells.Item(275, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(276, 9).Value = 'sv'
$ws.Cells.Item(276, 10).Value = 'Statement-opinion'
$ws.Cells.Item(283, 9).Value = 'sd'
$ws.Cells.Item(283, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(285, 9).Value = '%'
$ws.Cells.Item(285, 10).Value = 'Uninterpretable'
$ws.Cells.Item(305, 9).Value = 'sv'
$ws.Cells.Item(305, 10).Value = 'Statement-opinion'
$ws.Cells.Item(313, 9).Value = 'ba'
$ws.Cells.Item(313, 10).Value = 'Appreciation'
$ws.Cells.Item(317, 9).Value = 'sd'
$ws.Cells.Item(317, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(319, 9).Value = 'sd'
$ws.Cells.Item(319, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(320, 9).Value = 'sv'
$ws.Cells.Item(320, 10).Value = 'Statement-opinion'
$ws.Cells.Item(323, 9).Value = 'sd'
$ws.Cells.Item(323, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(324, 9).Value = 'sd'
$ws.Cells.Item(324, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(325, 9).Value = 'sd'
$ws.Cells.Item(325, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(326, 9).Value = 'sd'
$ws.Cells.Item(326, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(327, 9).Value = 'ba'
$ws.Cells.Item(327, 10).Value = 'Appreciation'
$ws.Cells.Item(337, 9).Value = 'sd'
$ws.Cells.Item(337, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(344, 9).Value = 'sd'
$ws.Cells.Item(344, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(346, 9).Value = 'sd'
$ws.Cells.Item(346, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(348, 9).Value = 'ba'
$ws.Cells.Item(348, 10).Value = 'Appreciation'
$ws.Cells.Item(388, 9).Value = 'b'
$ws.Cells.Item(388, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(395, 9).Value = 'b'
$ws.Cells.Item(395, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(396, 9).Value = 'sv'
$ws.Cells.Item(396, 10).Value = 'Statement-opinion'
$ws.Cells.Item(400, 9).Value = 'ba'
$ws.Cells.Item(400, 10).Value = 'Appreciation'
$ws.Cells.Item(403, 9).Value = 'sv'
$ws.Cells.Item(403, 10).Value = 'Statement-opinion'
$ws.Cells.Item(416, 9).Value = 'aa'
$ws.Cells.Item(416, 10).Value = 'Agree/Accept'
$ws.Cells.Item(418, 9).Value = 'aa'
$ws.Cells.Item(418, 10).Value = 'Agree/Accept'
$ws.Cells.Item(419, 9).Value = 'aa'
$ws.Cells.Item(419, 10).Value = 'Agree/Accept'
$ws.Cells.Item(440, 9).Value = 'b'
$ws.Cells.Item(440, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(450, 9).Value = 'ba'
$ws.Cells.Item(450, 10).Value = 'Appreciation'
$ws.Cells.Item(453, 9).Value = '%'
$ws.Cells.Item(453, 10).Value = 'Uninterpretable'
$ws.Cells.Item(458, 9).Value = 'sv'
$ws.Cells.Item(458, 10).Value = 'Statement-opinion'
$ws.Cells.Item(463, 9).Value = 'sv'
$ws.Cells.Item(463, 10).Value = 'Statement-opinion'
$ws.Cells.Item(472, 9).Value = 'sv'
$ws.Cells.Item(472, 10).Value = 'Statement-opinion'
